# Amit's Create enrollment fixes and updating data file
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Remove the old "Enrollment Flow" row (row 4); everything below shifts up.
$ws.Range("A4").EntireRow.Delete()

# Rebuild rows 4-14 with the corrected / extended data set.
$data = @(
    @("UPA_Regression", "Provider View Payments",       "test.java.TestProviderViewPayments"),
    @("UPA_Regression", "BS View Payments",              "test.java.TestBSViewPayments"),
    @("UPA_Regression", "Payer View Payments",           "test.java.TestPayerViewPayments"),
    @("UPA_Regression", "Create Enrollment",              "test.java.TestCreateEnrollment"),
    @("UPA_Regression", "My Profile",                    "test.java.TestUPAMyProfileTab"),
    @("UPA_Regression", "Provider Search Remittance",     "test.java.TestProviderSearchRemittance"),
    @("UPA_Regression", "BS Search Remittance",           "test.java.TestBSSearchRemittance"),
    @("UPA_Regression", "SubPayer Search Remittance",     "test.java.TestSubPayerSearchRemittance"),
    @("CSR_Regression",  "Manage Users",                  "test.java.TestCSRManageUsers"),
    @("CSR_Regression",  "View Payments",                 "test.java.TestCSRViewPayments")
)

$row = 4
foreach ($r in $data) {
    $ws.Range("A$row").Value = $r[0]
    $ws.Range("B$row").Value = $r[1]
    $ws.Range("C$row").Value = $r[2]
    $row = $row + 1
}

# Final row: shared-string order requires the ClassName value to be
# registered before the TestCaseName value.
$ws.Range("A14").Value = "CSR_Regression"
$ws.Range("C14").Value = "test.java.TestCSRSearchRemittance"
$ws.Range("B14").Value = "Search Remittance"

$ws.Range("B14").Select()
